$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E28: "PRESTAMOS" -> "PLANTILLA DE APORTES"
$ws.Cells.Item(28, 5).Value = "PLANTILLA DE APORTES"

# Append new rows 34-47
# Row 34
$ws.Cells.Item(34, 1).Value = "'1"
$ws.Cells.Item(34, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(34, 3).Value = "'1728220441001"
$ws.Cells.Item(34, 4).Value = "'1728220441"
$ws.Cells.Item(34, 5).Value = "PLANTILLA DE APORTES"
$ws.Cells.Item(34, 6).Value = 24.58
$ws.Cells.Item(34, 7).Value = 125.36
$ws.Cells.Item(34, 8).Value = "'162.00"
$ws.Cells.Item(34, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(34, 10).Value = "Mgs. Tibanta Tibanta"

# Row 35
$ws.Cells.Item(35, 1).Value = "'2"
$ws.Cells.Item(35, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(35, 3).Value = "'1728220441001"
$ws.Cells.Item(35, 4).Value = "'1728220441"
$ws.Cells.Item(35, 5).Value = "PLANTILLA DE APORTES"
$ws.Cells.Item(35, 6).Value = 152.36
$ws.Cells.Item(35, 7).Value = 258.58
$ws.Cells.Item(35, 8).Value = "'336.00"
$ws.Cells.Item(35, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(35, 10).Value = "Mgs. Tibanta Tibanta"

# Row 36
$ws.Cells.Item(36, 1).Value = "'3"
$ws.Cells.Item(36, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(36, 3).Value = "'1728220441001"
$ws.Cells.Item(36, 4).Value = "'1728220441"
$ws.Cells.Item(36, 5).Value = "PLANTILLA DE APORTES"
$ws.Cells.Item(36, 6).Value = 123.45
$ws.Cells.Item(36, 7).Value = 159.26
$ws.Cells.Item(36, 8).Value = "'207.00"
$ws.Cells.Item(36, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(36, 10).Value = "Mgs. Tibanta Tibanta"

# Row 37
$ws.Cells.Item(37, 1).Value = "'4"
$ws.Cells.Item(37, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(37, 3).Value = "'1728220441001"
$ws.Cells.Item(37, 4).Value = "'1728220441"
$ws.Cells.Item(37, 5).Value = "PLANILLA DE PRESTAMOS"
$ws.Cells.Item(37, 6).Value = 158.47
$ws.Cells.Item(37, 7).Value = 356.48
$ws.Cells.Item(37, 8).Value = "'463.00"
$ws.Cells.Item(37, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(37, 10).Value = "Mgs. Tibanta Tibanta"

# Row 38
$ws.Cells.Item(38, 1).Value = "'5"
$ws.Cells.Item(38, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(38, 3).Value = "'1728220441001"
$ws.Cells.Item(38, 4).Value = "'1728220441"
$ws.Cells.Item(38, 5).Value = "PLANILLA DE PRESTAMOS"
$ws.Cells.Item(38, 6).Value = 478.25
$ws.Cells.Item(38, 7).Value = 258.69
$ws.Cells.Item(38, 8).Value = "'336.00"
$ws.Cells.Item(38, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(38, 10).Value = "Mgs. Tibanta Tibanta"

# Row 39
$ws.Cells.Item(39, 1).Value = "'6"
$ws.Cells.Item(39, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(39, 3).Value = "'1728220441001"
$ws.Cells.Item(39, 4).Value = "'1728220441"
$ws.Cells.Item(39, 5).Value = "PLANILLA DE PRESTAMOS"
$ws.Cells.Item(39, 6).Value = "'300.25"
$ws.Cells.Item(39, 7).Value = 425.12
$ws.Cells.Item(39, 8).Value = "'552.00"
$ws.Cells.Item(39, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(39, 10).Value = "Mgs. Tibanta Tibanta"

# Row 40
$ws.Cells.Item(40, 1).Value = "'7"
$ws.Cells.Item(40, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(40, 3).Value = "'1728220441001"
$ws.Cells.Item(40, 4).Value = "'1728220441"
$ws.Cells.Item(40, 5).Value = "PLANILLA DE RESPONSABILIDAD PATRONAL"
$ws.Cells.Item(40, 6).Value = 456.25
$ws.Cells.Item(40, 7).Value = 789.125
$ws.Cells.Item(40, 8).Value = "'1025.00"
$ws.Cells.Item(40, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(40, 10).Value = "Mgs. Tibanta Tibanta"

# Row 41
$ws.Cells.Item(41, 1).Value = "'8"
$ws.Cells.Item(41, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(41, 3).Value = "'1728220441001"
$ws.Cells.Item(41, 4).Value = "'1728220441"
$ws.Cells.Item(41, 5).Value = "PLANILLA DE RESPONSABILIDAD PATRONAL"
$ws.Cells.Item(41, 6).Value = 456.258
$ws.Cells.Item(41, 7).Value = 500.12
$ws.Cells.Item(41, 8).Value = "'650.00"
$ws.Cells.Item(41, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(41, 10).Value = "Mgs. Tibanta Tibanta"

# Row 42
$ws.Cells.Item(42, 1).Value = "'9"
$ws.Cells.Item(42, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(42, 3).Value = "'1728220441001"
$ws.Cells.Item(42, 4).Value = "'1728220441"
$ws.Cells.Item(42, 5).Value = "PLANILLA DE FONDOS DE RESERVA"
$ws.Cells.Item(42, 6).Value = 45.12
$ws.Cells.Item(42, 7).Value = 100.25
$ws.Cells.Item(42, 8).Value = "'130.00"
$ws.Cells.Item(42, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(42, 10).Value = "Mgs. Tibanta Tibanta"

# Row 43
$ws.Cells.Item(43, 1).Value = "'10"
$ws.Cells.Item(43, 2).Value = "ALEXANDER TIBANTA MIRANDA FRANCISCO"
$ws.Cells.Item(43, 3).Value = "'1728220441001"
$ws.Cells.Item(43, 4).Value = "'1728220441"
$ws.Cells.Item(43, 5).Value = "PLANILLA DE FONDOS DE RESERVA"
$ws.Cells.Item(43, 6).Value = 66.58
$ws.Cells.Item(43, 7).Value = 75.12
$ws.Cells.Item(43, 8).Value = "'97.00"
$ws.Cells.Item(43, 9).Value = "Dr. Jorge Gonzalo Atiencia Gálvez"
$ws.Cells.Item(43, 10).Value = "Mgs. Tibanta Tibanta"

# Row 44
$ws.Cells.Item(44, 1).Value = "'4578"
$ws.Cells.Item(44, 2).Value = "BBBB"
$ws.Cells.Item(44, 3).Value = "'12457845001"
$ws.Cells.Item(44, 4).Value = "'12457855"
$ws.Cells.Item(44, 5).Value = "PLANILLA DE FONDOS DE RESERVA"
$ws.Cells.Item(44, 6).Value = 45.25
$ws.Cells.Item(44, 7).Value = 123.52
$ws.Cells.Item(44, 8).Value = "'160.00"
$ws.Cells.Item(44, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(44, 10).Value = "AAAAA"

# Row 45
$ws.Cells.Item(45, 1).Value = "'4579"
$ws.Cells.Item(45, 2).Value = "BBBB"
$ws.Cells.Item(45, 3).Value = "'12457845001"
$ws.Cells.Item(45, 4).Value = "'12457855"
$ws.Cells.Item(45, 5).Value = "PLANILLA DE FONDOS DE RESERVA"
$ws.Cells.Item(45, 6).Value = 45.25
$ws.Cells.Item(45, 7).Value = 123.52
$ws.Cells.Item(45, 8).Value = "'160.00"
$ws.Cells.Item(45, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(45, 10).Value = "AAAAA"

# Row 46
$ws.Cells.Item(46, 1).Value = "'4580"
$ws.Cells.Item(46, 2).Value = "BBBB"
$ws.Cells.Item(46, 3).Value = "'12457845001"
$ws.Cells.Item(46, 4).Value = "'12457855"
$ws.Cells.Item(46, 5).Value = "PLANILLA DE FONDOS DE RESERVA"
$ws.Cells.Item(46, 6).Value = 45.25
$ws.Cells.Item(46, 7).Value = 123.52
$ws.Cells.Item(46, 8).Value = "'160.00"
$ws.Cells.Item(46, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(46, 10).Value = "AAAAA"

# Row 47
$ws.Cells.Item(47, 1).Value = "'124578"
$ws.Cells.Item(47, 2).Value = "BBBB"
$ws.Cells.Item(47, 3).Value = "'12457845001"
$ws.Cells.Item(47, 4).Value = "'12457855"
$ws.Cells.Item(47, 5).Value = "PLANTILLA DE APORTES"
$ws.Cells.Item(47, 6).Value = 12.25
$ws.Cells.Item(47, 7).Value = 123.25
$ws.Cells.Item(47, 8).Value = "'160.00"
$ws.Cells.Item(47, 9).Value = "Dr. Christian Santiago Izurieta Cruz"
$ws.Cells.Item(47, 10).Value = "AAAAA"
